# Weekly data refresh: insert a new price record for this week at row 138
# (pushing all subsequent rows down by one) on the Berenjena / Terminal La
# Palmera de La Serena sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 138 - this shifts old rows 138:215
# down to 139:216 and grows the sheet's used range to A1:R216 automatically.
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with this week's record.
$ws.Range("A138").Value = 8
$ws.Range("B138").Value = "Terminal La Palmera de La Serena"
$ws.Range("C138").Value = "Coquimbo"
$ws.Range("D138").Value = 45001
$ws.Range("E138").Value = 4
$ws.Range("F138").Value = 100112001
$ws.Range("G138").Value = "Berenjena"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 440
$ws.Range("K138").Value = 11800
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = 11900
$ws.Range("N138").Value = "$/caja 50 unidades"
$ws.Range("O138").Value = "Región de Arica y Parinacota"
$ws.Range("P138").Value = 238
$ws.Range("Q138").Value = 50
$ws.Range("R138").Value = "Hortaliza"
